$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 46064

# Row 3
$ws.Range("C3").Value = 46064

# Row 4
$ws.Range("C4").Value = 46064

# Row 5
$ws.Range("C5").Value = 46064

# Row 6
$ws.Range("C6").Value = 46064

# Row 7
$ws.Range("C7").Value = 46064

# Row 8
$ws.Range("C8").Value = 46064

# Row 9
$ws.Range("C9").Value = 46064

# Row 10
$ws.Range("C10").Value = 46064

# Row 11
$ws.Range("A11").Value = "A 44069-2021"
$ws.Range("B11").Value = 44434
$ws.Range("C11").Value = 46064
$ws.Range("G11").Value = 1.5

# Row 12
$ws.Range("A12").Value = "A 33953-2021"
$ws.Range("B12").Value = 44378
$ws.Range("C12").Value = 46064
$ws.Range("G12").Value = 5

# Row 13
$ws.Range("C13").Value = 46064

# Row 14
$ws.Range("A14").Value = "A 21264-2022"
$ws.Range("B14").Value = 44705
$ws.Range("C14").Value = 46064
$ws.Range("F14").Value = "Allmännings- och besparingsskogar"
$ws.Range("G14").Value = 2.4

# Row 15
$ws.Range("A15").Value = "A 21972-2023"
$ws.Range("B15").Value = 45068.66232638889
$ws.Range("C15").Value = 46064
$ws.Range("G15").Value = 1.5

# Row 16
$ws.Range("A16").Value = "A 41546-2025"
$ws.Range("B16").Value = 45901.57927083333
$ws.Range("C16").Value = 46064
$ws.Range("F16").Value = "Allmännings- och besparingsskogar"
$ws.Range("G16").Value = 4.3

# Row 17
$ws.Range("A17").Value = "A 41550-2025"
$ws.Range("B17").Value = 45901.58652777778
$ws.Range("C17").Value = 46064
$ws.Range("F17").Value = "Allmännings- och besparingsskogar"
$ws.Range("G17").Value = 2.3

# Row 18
$ws.Range("A18").Value = "A 62831-2023"
$ws.Range("B18").Value = 45270
$ws.Range("C18").Value = 46064
$ws.Range("G18").Value = 1.1

# Row 19
$ws.Range("A19").Value = "A 42991-2025"
$ws.Range("B19").Value = 45909.45190972222
$ws.Range("C19").Value = 46064
$ws.Range("G19").Value = 7.2

# Row 20
$ws.Range("A20").Value = "A 42994-2025"
$ws.Range("B20").Value = 45909.453518518516
$ws.Range("C20").Value = 46064
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = 7.9

# Row 21
$ws.Range("A21").Value = "A 10263-2024"
$ws.Range("B21").Value = 45365.43090277778
$ws.Range("C21").Value = 46064
$ws.Range("F21").Value = "Kyrkan"
$ws.Range("G21").Value = 1.4

# Row 22
$ws.Range("A22").Value = "A 30743-2021"
$ws.Range("B22").Value = 44365
$ws.Range("C22").Value = 46064
$ws.Range("F22").Value = ""
$ws.Range("G22").Value = 3

# Row 23
$ws.Range("A23").Value = "A 43448-2025"
$ws.Range("B23").Value = 45911.45209490741
$ws.Range("C23").Value = 46064
$ws.Range("G23").Value = 1.2

# Row 24
$ws.Range("A24").Value = "A 55562-2022"
$ws.Range("B24").Value = 44888
$ws.Range("C24").Value = 46064
$ws.Range("G24").Value = 0.8

# Row 25
$ws.Range("A25").Value = "A 44192-2025"
$ws.Range("B25").Value = 45915.61556712963
$ws.Range("C25").Value = 46064
$ws.Range("G25").Value = 0.8

# Row 26
$ws.Range("C26").Value = 46064

# Row 27
$ws.Range("A27").Value = "A 12077-2022"
$ws.Range("B27").Value = 44636.47484953704
$ws.Range("C27").Value = 46064
$ws.Range("G27").Value = 2.1

# Row 28
$ws.Range("A28").Value = "A 27365-2025"
$ws.Range("B28").Value = 45812.64355324074
$ws.Range("C28").Value = 46064
$ws.Range("G28").Value = 11.9

# Row 29
$ws.Range("A29").Value = "A 37407-2023"
$ws.Range("B29").Value = 45156.60152777778
$ws.Range("C29").Value = 46064
$ws.Range("G29").Value = 3.3

# Row 30
$ws.Range("A30").Value = "A 30174-2021"
$ws.Range("B30").Value = 44363
$ws.Range("C30").Value = 46064
$ws.Range("G30").Value = 1.8

# Row 31
$ws.Range("A31").Value = "A 33799-2025"
$ws.Range("B31").Value = 45842.440567129626
$ws.Range("C31").Value = 46064
$ws.Range("F31").Value = "Allmännings- och besparingsskogar"
$ws.Range("G31").Value = 2.5

# Row 32
$ws.Range("A32").Value = "A 37072-2025"
$ws.Range("B32").Value = 45875.40431712963
$ws.Range("C32").Value = 46064
$ws.Range("G32").Value = 1.5

# Row 33
$ws.Range("A33").Value = "A 37076-2025"
$ws.Range("B33").Value = 45875.41342592592
$ws.Range("C33").Value = 46064
$ws.Range("G33").Value = 2.4

# Row 34
$ws.Range("A34").Value = "A 13510-2025"
$ws.Range("B34").Value = 45736.471030092594
$ws.Range("C34").Value = 46064
$ws.Range("G34").Value = 2.5

# Row 35
$ws.Range("A35").Value = "A 14149-2022"
$ws.Range("B35").Value = 44651
$ws.Range("C35").Value = 46064
$ws.Range("F35").Value = "Allmännings- och besparingsskogar"
$ws.Range("G35").Value = 3.8

# Row 36
$ws.Range("A36").Value = "A 53343-2024"
$ws.Range("B36").Value = 45614.43885416666
$ws.Range("C36").Value = 46064
$ws.Range("G36").Value = 0.9

# Row 37
$ws.Range("A37").Value = "A 8848-2025"
$ws.Range("B37").Value = 45713.34292824074
$ws.Range("C37").Value = 46064
$ws.Range("G37").Value = 2

# Row 38
$ws.Range("A38").Value = "A 8436-2023"
$ws.Range("B38").Value = 44977
$ws.Range("C38").Value = 46064
$ws.Range("F38").Value = "Kyrkan"
$ws.Range("G38").Value = 4

# Row 39
$ws.Range("A39").Value = "A 58109-2025"
$ws.Range("B39").Value = 45982.597650462965
$ws.Range("C39").Value = 46064
$ws.Range("G39").Value = 2.5

# Row 40
$ws.Range("A40").Value = "A 58111-2025"
$ws.Range("B40").Value = 45982.59920138889
$ws.Range("C40").Value = 46064
$ws.Range("G40").Value = 0.6

# Row 41
$ws.Range("A41").Value = "A 58125-2025"
$ws.Range("B41").Value = 45982.615069444444
$ws.Range("C41").Value = 46064
$ws.Range("G41").Value = 0.8

# Row 42
$ws.Range("A42").Value = "A 54207-2025"
$ws.Range("B42").Value = 45964
$ws.Range("C42").Value = 46064
$ws.Range("G42").Value = 2.1

# Row 43
$ws.Range("A43").Value = "A 54203-2025"
$ws.Range("B43").Value = 45964
$ws.Range("C43").Value = 46064
$ws.Range("G43").Value = 10.3

# Row 44
$ws.Range("A44").Value = "A 7245-2025"
$ws.Range("B44").Value = 45702
$ws.Range("C44").Value = 46064
$ws.Range("G44").Value = 4

# Row 45
$ws.Range("A45").Value = "A 1621-2026"
$ws.Range("B45").Value = 46034.47645833333
$ws.Range("C45").Value = 46064
$ws.Range("G45").Value = 1.4

# Row 46
$ws.Range("A46").Value = "A 1622-2026"
$ws.Range("B46").Value = 46034.47929398148
$ws.Range("C46").Value = 46064
$ws.Range("F46").Value = ""
$ws.Range("G46").Value = 1.4

# Row 47
$ws.Range("A47").Value = "A 34984-2024"
$ws.Range("B47").Value = 45527
$ws.Range("C47").Value = 46064
$ws.Range("F47").Value = ""
$ws.Range("G47").Value = 4.1

# Row 48
$ws.Range("A48").Value = "A 62433-2025"
$ws.Range("B48").Value = 46007
$ws.Range("C48").Value = 46064
$ws.Range("G48").Value = 2

# Row 49
$ws.Range("A49").Value = "A 16762-2022"
$ws.Range("B49").Value = 44673.47876157407
$ws.Range("C49").Value = 46064
$ws.Range("G49").Value = 4.2

# Row 50
$ws.Range("A50").Value = "A 22072-2023"
$ws.Range("B50").Value = 45069
$ws.Range("C50").Value = 46064
$ws.Range("F50").Value = ""
$ws.Range("G50").Value = 3.5

# Row 51
$ws.Range("A51").Value = "A 30766-2022"
$ws.Range("B51").Value = 44764
$ws.Range("C51").Value = 46064
$ws.Range("G51").Value = 0.6

# Row 52
$ws.Range("A52").Value = "A 35036-2024"
$ws.Range("B52").Value = 45527
$ws.Range("C52").Value = 46064
$ws.Range("G52").Value = 1.7

